$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force account id column to be stored as text so leading zeros survive,
# then drop the format override again so the cells keep plain default style.
$ws.Range("B1:B5").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "resourceName"
$ws.Range("B1").Value = "accountIds"
$ws.Range("C1").Value = "regions"

# Row 2
$ws.Range("A2").Value = "testResource"
$ws.Range("B2").Value = "060087218145"
$ws.Range("C2").Value = "Mumbai"

# Row 3
$ws.Range("A3").Value = "resourceSingapore"
$ws.Range("B3").Value = "060087218145"
$ws.Range("C3").Value = "Singapore"

# Row 4
$ws.Range("A4").Value = "ShubhamTest"
$ws.Range("B4").Value = "060087218145"
$ws.Range("C4").Value = "Singapore"

# Row 5
$ws.Range("A5").Value = "resourceMumbai"
$ws.Range("B5").Value = "060087218145"
$ws.Range("C5").Value = "Mumbai"

# Remove the temporary text-format override so the cells end up with the
# default (unstyled) appearance, matching the original sheet's formatting.
$ws.Range("B1:B5").ClearFormats()
